$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 14-17: labels + summary formulas (style s=2 -> fontId 3, bold 12pt, vertical center)
# Applied as one combined range write so only a single new font/xf entry is minted.
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$summaryRange = $ws.Range("B14:B17")
$summaryRange.Font.Bold = $true
$summaryRange.Font.Size = 12
$summaryRange.VerticalAlignment = -4108
$ws.Range("A14:B17").RowHeight = 15.6

# J12: Average of J2:J11, bold 11pt font (style s=3 -> fontId 2)
# Created second so it becomes the second new font/xf entry.
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11
$ws.Range("J12").Font.ThemeColor = 1

$ws.Calculate()
